$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text, preserving exact formatting (avoids numeric coercion)
function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextCell "D2" "27.466.92"
Set-TextCell "D3" "1.863.88"
Set-TextCell "E3" "  +0.93%  "
Set-TextCell "D4" "1.011"
Set-TextCell "E4" "  -0.23%  "
Set-TextCell "D5" "310.98"
Set-TextCell "E5" "  +0.33%  "
Set-TextCell "E6" "  -0.16%  "
Set-TextCell "D7" "0.4781"
Set-TextCell "E7" "  +0.04%  "
Set-TextCell "D8" "0.3811"
Set-TextCell "E8" "  +3.61%  "
Set-TextCell "D9" "0.07328"
Set-TextCell "E9" "  +1.34%  "
Set-TextCell "D10" "0.9366"
Set-TextCell "E10" "  +0.84%  "
Set-TextCell "D11" "20.75"
Set-TextCell "E11" "  +5.37%  "
Set-TextCell "E12" "  +0.84%  "
Set-TextCell "D13" "1.879.52"
Set-TextCell "E13" "  +2.61%  "
Set-TextCell "E14" "  +1.91%  "
Set-TextCell "D15" "6.563"
Set-TextCell "D16" "90.61"
Set-TextCell "E16" "  +2.09%  "
Set-TextCell "E17" "  -0.22%  "
Set-TextCell "D18" "0.000008807"
Set-TextCell "E18" "  +1.91%  "
Set-TextCell "E19" "  -0.23%  "
Set-TextCell "D20" "27.463.82"
Set-TextCell "E20" "  +1.70%  "
Set-TextCell "E21" "  +1.54%  "
Set-TextCell "D22" "5.119"
Set-TextCell "E22" "  +1.15%  "
Set-TextCell "E23" "  +0.58%  "
Set-TextCell "D24" "1.941"
Set-TextCell "E24" "  +1.01%  "
Set-TextCell "D25" "154.93"
Set-TextCell "E25" "  +1.38%  "
Set-TextCell "D26" "18.50"
Set-TextCell "E26" "  +1.64%  "
Set-TextCell "D27" "2.021"
Set-TextCell "E27" "  +1.06%  "
Set-TextCell "D28" "115.57"
Set-TextCell "E28" "  +1.12%  "
Set-TextCell "D29" "4.951"
Set-TextCell "E29" "  -0.40%  "
Set-TextCell "D30" "0.08897"
Set-TextCell "E30" "  +0.05%  "
Set-TextCell "D31" "3.321"
Set-TextCell "E31" "  -0.10%  "
Set-TextCell "E32" "  +3.67%  "
Set-TextCell "D33" "0.7599"
Set-TextCell "E33" "  +2.30%  "
Set-TextCell "D34" "4.610"
Set-TextCell "E34" "  +2.38%  "
Set-TextCell "D35" "2.732"
Set-TextCell "E35" "  -0.75%  "
Set-TextCell "D36" "0.02059"
Set-TextCell "E36" "  +5.02%  "
Set-TextCell "E37" "  +0.08%  "
Set-TextCell "D38" "0.5591"
Set-TextCell "E38" "  +7.15%  "
Set-TextCell "D39" "0.05281"
Set-TextCell "E39" "  +0.13%  "
Set-TextCell "E40" "  +0.50%  "
Set-TextCell "D41" "7.066"
Set-TextCell "E41" "  +1.04%  "
Set-TextCell "D42" "8.679"
Set-TextCell "E42" "  +5.57%  "
Set-TextCell "D43" "0.1529"
Set-TextCell "E43" "  +1.03%  "
Set-TextCell "D44" "0.4913"
Set-TextCell "E44" "  +3.24%  "
Set-TextCell "D45" "10.73"
Set-TextCell "E45" "  +0.62%  "
Set-TextCell "D46" "1.011"
Set-TextCell "D49" "67.44"
Set-TextCell "E49" "  +2.91%  "
Set-TextCell "D50" "0.06081"
Set-TextCell "E50" "  +0.30%  "
Set-TextCell "D51" "0.9154"
Set-TextCell "E51" "  +3.11%  "

# Row 47/48: Quant and NEARProtocol swapped position with updated figures
Set-TextCell "B47" "Quant"
Set-TextCell "C47" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D47" "103.13"
Set-TextCell "E47" "  +1.36%  "
Set-TextCell "B48" "NEARProtocol"
Set-TextCell "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D48" "1.658"
Set-TextCell "E48" "  +3.04%  "
